$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (Camote, Femacal de La Calera / Zapallo) needs to be
# inserted as the new first entry of this block, pushing the existing rows
# 341-355 down to 342-356 (dimension grows from R355 to R356).
$ws.Rows.Item(341).Insert()

# Populate the newly inserted row 341 with the new record's data.
$ws.Cells.Item(341, 1).Value = 3
$ws.Cells.Item(341, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(341, 3).Value = "Coquimbo"
$ws.Cells.Item(341, 4).Value = 44509
$ws.Cells.Item(341, 5).Value = 5
$ws.Cells.Item(341, 6).Value = 100112045
$ws.Cells.Item(341, 7).Value = "Zapallo"
$ws.Cells.Item(341, 8).Value = "Camote"
$ws.Cells.Item(341, 9).Value = "1a nueva(o)"
$ws.Cells.Item(341, 10).Value = 160
$ws.Cells.Item(341, 11).Value = 500
$ws.Cells.Item(341, 12).Value = 500
$ws.Cells.Item(341, 13).Value = 500
$ws.Cells.Item(341, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(341, 15).Value = "Provincia de Talca"
$ws.Cells.Item(341, 16).Value = 500
$ws.Cells.Item(341, 17).Value = 1
$ws.Cells.Item(341, 18).Value = "Hortaliza"
